$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Clear out the old "Venda" schema table (A1:C7) ---
$ws.Range("A1:G11").ClearContents()

# --- Left table: Order ---
$ws.Range("A1").Value = "Order"

$ws.Range("A2").Value = "Id"
$ws.Range("B2").Value = "varchar(32)"
$ws.Range("C2").Value = "primary key"

$ws.Range("A3").Value = "CustomerId"
$ws.Range("B3").Value = "varchar(32)"

$ws.Range("A4").Value = "Status"
$ws.Range("B4").Value = "TINYINT"

$ws.Range("A5").Value = "Total"
$ws.Range("B5").Value = "decimal(14,2)"
$ws.Range("C5").Value = "notnull"

$ws.Range("A6").Value = "Created"
$ws.Range("B6").Value = "date"
$ws.Range("C6").Value = "notnull"

$ws.Range("A7").Value = "Updated"
$ws.Range("B7").Value = "date"
$ws.Range("C7").Value = "notnull"

$ws.Range("A8").Value = "PaymentForm"
$ws.Range("B8").Value = "TINYINT"

# --- Right table: OrderItem ---
$ws.Range("E1").Value = "OrderItem"

$ws.Range("E2").Value = "OrderId"
$ws.Range("F2").Value = "varchar(32)"
$ws.Range("G2").Value = "primary key"

$ws.Range("E3").Value = "ProductId"
$ws.Range("F3").Value = "varchar(32)"
$ws.Range("G3").Value = "primary key"

$ws.Range("E4").Value = "Count"
$ws.Range("F4").Value = "integer"
$ws.Range("G4").Value = "notnull"

$ws.Range("E5").Value = "UnitValue"
$ws.Range("F5").Value = "decimal(14,2)"
$ws.Range("G5").Value = "notnull"

$ws.Range("E6").Value = "Total"
$ws.Range("F6").Value = "decimal(14,2)"
$ws.Range("G6").Value = "notnull"

$ws.Range("E7").Value = "ProductName"
$ws.Range("F7").Value = "varchar(100)"

# --- Database name ---
$ws.Range("A11").Value = "RLSalesDB"

# --- Column widths for the new right-hand table ---
# (target stored widths are 17 / 14.7109375 / 12.85546875 "characters";
#  the host's ColumnWidth setter quantizes internally to 1/6-character
#  steps, so we pick the input that lands on the closest reachable step)
$ws.Columns.Item(5).ColumnWidth = 16.1
$ws.Columns.Item(6).ColumnWidth = 13.8
$ws.Columns.Item(7).ColumnWidth = 12

# --- Selection moves to A12, matching the post-edit cursor position ---
$ws.Range("A12").Select() | Out-Null
